$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.39646633333334
$ws.Range("H2").Value = 64.18939900000001
$ws.Range("I2").Value = 0.0721325008796955
$ws.Range("J2").Value = 0.0721325008796955
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 39.68128022858723
$ws.Range("R2").Value = 357.1315220572851
$ws.Range("S2").Value = 0.001191175667189359
$ws.Range("T2").Value = 0.001191175667189359
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.39646633333334
$ws.Range("H3").Value = 64.18939900000001
$ws.Range("I3").Value = 0.0721325008796955
$ws.Range("J3").Value = 0.0721325008796955
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 1988.025374602781
$ws.Range("R3").Value = 17892.22837142503
$ws.Range("S3").Value = 0.05967769785501584
$ws.Range("T3").Value = 0.05967769785501585
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.39646633333334
$ws.Range("H4").Value = 64.18939900000001
$ws.Range("I4").Value = 0.0721325008796955
$ws.Range("J4").Value = 0.0721325008796955
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 3.469137465421334
$ws.Range("R4").Value = 31.22223718879201
$ws.Range("S4").Value = 0.0001041385789757818
$ws.Range("T4").Value = 0.0001041385789757818
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.39646633333334
$ws.Range("H5").Value = 64.18939900000001
$ws.Range("I5").Value = 0.0721325008796955
$ws.Range("J5").Value = 0.0721325008796955
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 371.752726004614
$ws.Range("R5").Value = 3345.774534041525
$ws.Range("S5").Value = 0.01115948877851451
$ws.Range("T5").Value = 0.01115948877851451
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 130.868154
$ws.Range("H6").Value = 392.604462
$ws.Range("I6").Value = 0.441187207572817
$ws.Range("J6").Value = 0.441187207572817
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 242.70437047737
$ws.Range("R6").Value = 2184.33933429633
$ws.Range("S6").Value = 0.007285640452317822
$ws.Range("T6").Value = 0.007285640452317824
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 130.868154
$ws.Range("H7").Value = 392.604462
$ws.Range("I7").Value = 0.441187207572817
$ws.Range("J7").Value = 0.441187207572817
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 12159.44758476822
$ws.Range("R7").Value = 109435.0282629139
$ws.Range("S7").Value = 0.365009344607932
$ws.Range("T7").Value = 0.3650093446079321
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 130.868154
$ws.Range("H8").Value = 392.604462
$ws.Range("I8").Value = 0.441187207572817
$ws.Range("J8").Value = 0.441187207572817
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 21.218439016944
$ws.Range("R8").Value = 190.965951152496
$ws.Range("S8").Value = 0.0006369473995578508
$ws.Range("T8").Value = 0.0006369473995578509
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 130.868154
$ws.Range("H9").Value = 392.604462
$ws.Range("I9").Value = 0.441187207572817
$ws.Range("J9").Value = 0.441187207572817
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 2273.76765110505
$ws.Range("R9").Value = 20463.90885994545
$ws.Range("S9").Value = 0.06825527511300934
$ws.Range("T9").Value = 0.06825527511300934
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 75.02619166666666
$ws.Range("H10").Value = 225.078575
$ws.Range("I10").Value = 0.2529308696158396
$ws.Range("J10").Value = 0.2529308696158397
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 139.1414493229028
$ws.Range("R10").Value = 1252.273043906125
$ws.Range("S10").Value = 0.004176828665207709
$ws.Range("T10").Value = 0.004176828665207711
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 75.02619166666666
$ws.Range("H11").Value = 225.078575
$ws.Range("I11").Value = 0.2529308696158396
$ws.Range("J11").Value = 0.2529308696158397
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 6970.962890296497
$ws.Range("R11").Value = 62738.66601266847
$ws.Range("S11").Value = 0.2092584040627568
$ws.Range("T11").Value = 0.2092584040627569
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 75.02619166666666
$ws.Range("H12").Value = 225.078575
$ws.Range("I12").Value = 0.2529308696158396
$ws.Range("J12").Value = 0.2529308696158397
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 12.16444661206667
$ws.Range("R12").Value = 109.4800195086
$ws.Range("S12").Value = 0.0003651594083065635
$ws.Range("T12").Value = 0.0003651594083065636
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 75.02619166666666
$ws.Range("H13").Value = 225.078575
$ws.Range("I13").Value = 0.2529308696158396
$ws.Range("J13").Value = 0.2529308696158397
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 1303.541941894236
$ws.Range("R13").Value = 11731.87747704813
$ws.Range("S13").Value = 0.03913047747956849
$ws.Range("T13").Value = 0.0391304774795685
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 69.33645133333333
$ws.Range("H14").Value = 208.009354
$ws.Range("I14").Value = 0.2337494219316478
$ws.Range("J14").Value = 0.2337494219316478
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 128.5894181100122
$ws.Range("R14").Value = 1157.30476299011
$ws.Range("S14").Value = 0.003860071676829027
$ws.Range("T14").Value = 0.003860071676829028
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 69.33645133333333
$ws.Range("H15").Value = 208.009354
$ws.Range("I15").Value = 0.2337494219316478
$ws.Range("J15").Value = 0.2337494219316478
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 6442.307925436916
$ws.Range("R15").Value = 57980.77132893224
$ws.Range("S15").Value = 0.1933889329455948
$ws.Range("T15").Value = 0.1933889329455948
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 69.33645133333333
$ws.Range("H16").Value = 208.009354
$ws.Range("I16").Value = 0.2337494219316478
$ws.Range("J16").Value = 0.2337494219316478
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 11.24193487338133
$ws.Range("R16").Value = 101.177413860432
$ws.Range("S16").Value = 0.0003374669162929901
$ws.Range("T16").Value = 0.0003374669162929902
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 69.33645133333333
$ws.Range("H17").Value = 208.009354
$ws.Range("I17").Value = 0.2337494219316478
$ws.Range("J17").Value = 0.2337494219316478
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 1204.685595887239
$ws.Range("R17").Value = 10842.17036298515
$ws.Range("S17").Value = 0.03616295039293096
$ws.Range("T17").Value = 0.03616295039293096
